$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G, rows 3 and 4
# both referenced the same shared string "2016-08-29 08:17:21" -> "2016-08-29 08:18:23"
$wsOverview.Range("G3").Value = "2016-08-29 08:18:23"
$wsOverview.Range("G4").Value = "2016-08-29 08:18:23"

# zh-cn sheet, row 3 (61cd30f7 file):
#  E3 "Priority": ht -> mt  (also shared with E4)
#  H3 "Correspond Handoff Datetime": 2016-08-29 08:17:15 -> 2016-08-29 08:18:18
#  K3 "Correspond Handback DateTime": 2016-08-29 08:17:42 -> 2016-08-29 08:18:43
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-29 08:18:18"
$wsZhCn.Range("K3").Value = "2016-08-29 08:18:43"

# de-de sheet, row 3 (61cd30f7 file):
#  E3 "Priority": ht -> mt (also shared with E4)
#  H3 "Correspond Handoff Datetime" shares the same string as Overview G3/G4
#  K3 "Correspond Handback DateTime": 2016-08-29 08:17:49 -> 2016-08-29 08:18:50
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-29 08:18:23"
$wsDeDe.Range("H4").Value = "2016-08-29 08:18:23"
$wsDeDe.Range("K3").Value = "2016-08-29 08:18:50"
